$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.436.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "'3.629.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("E4").Value = "  +43.12%  "
$ws.Range("D5").Value = "'1.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'226.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.17%  "
$ws.Range("D7").Value = "'641.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.28%  "
$ws.Range("D8").Value = "'0.424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").Value = "'1.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.02%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'3.628.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").Value = "'48.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.56%  "
$ws.Range("D13").Value = "'0.213"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("E14").Value = "  -9.39%  "
$ws.Range("D15").Value = "'6.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").Value = "'4.309.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "'96.175.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'21.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.70%  "
$ws.Range("D19").Value = "'8.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").Value = "'13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.42%  "
$ws.Range("D21").Value = "'3.631.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "'0.563"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.15%  "
$ws.Range("D23").Value = "'0.279"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +46.63%  "
$ws.Range("D24").Value = "'517.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("E25").Value = "  -7.13%  "
$ws.Range("D26").Value = "'120.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.55%  "
$ws.Range("D27").Value = "'0.0000201"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -12.46%  "
$ws.Range("D28").Value = "'6.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'3.815.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "'12.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("D31").Value = "'12.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "'3.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'0.621"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").Value = "'32.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").Value = "'0.179"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.08%  "
$ws.Range("D38").Value = "'1.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'8.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("D41").Value = "'586.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.66%  "
$ws.Range("D42").Value = "'7.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("D43").Value = "'0.500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.92%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0506"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.56%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'40.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("D47").Value = "'0.956"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").Value = "'1.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "'230.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.23%  "
$ws.Range("D50").Value = "'8.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'2.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.32%  "
